# Applies the updated Betfair Back/Lay odds for 2025-12-29 to Sheet1.
# Only numeric cell values change; no structural/style changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "AN2" = 7.4
    "F2"  = 1.64
    "H2"  = 5.5
    "S2"  = 2.7
    "T2"  = 1.76
    "U2"  = 2.24
    "I3"  = 5.1
    "K3"  = 3.3
    "Q3"  = 2.68
    "V3"  = 1.26
    "F4"  = 2.4
    "G4"  = 2.52
    "H4"  = 2.82
    "I4"  = 3
    "R4"  = 1.62
    "S4"  = 2.3
    "U4"  = 2.7
    "V4"  = 1.5
    "W4"  = 1.65
    "AF5" = 16.5
    "F5"  = 2.42
    "G5"  = 2.74
    "H5"  = 3.45
    "I5"  = 4.1
    "J5"  = 2.74
    "L5"  = 1.65
    "O5"  = 1.66
    "Q5"  = 3
    "T5"  = 2.28
    "V5"  = 1.33
    "W5"  = 1.57
    "X5"  = 7.6
    "Z5"  = 29
    "AD6" = 32
    "G6"  = 1.62
    "M6"  = 1.05
    "S6"  = 2.74
    "T6"  = 1.85
    "U6"  = 1.96
    "V6"  = 1.14
    "N7"  = 1.1
    "S7"  = 2.56
    "F8"  = 1.55
    "G8"  = 1.74
    "H8"  = 8
    "I8"  = 15
    "K8"  = 4.2
    "L8"  = 1.56
    "N8"  = 2.24
    "O8"  = 1.65
    "P8"  = 1.42
    "Q8"  = 2.72
    "R8"  = 1.1
    "S8"  = 4.5
    "V8"  = 1.1
    "W8"  = 2.34
    "AL9" = 55
    "AM9" = 240
    "AN9" = 14.5
    "G9"  = 1.68
    "N9"  = 3
    "O9"  = 1.48
    "P9"  = 1.66
    "Q9"  = 2.44
    "U9"  = 1.69
    "W9"  = 2.46
    "AC10" = 28
    "AG10" = 16.5
    "AI10" = 540
    "AL10" = 70
    "H10" = 25
    "J10" = 11
    "K10" = 12
    "L10" = 1.22
    "N10" = 7
    "O10" = 1.14
    "P10" = 2.96
    "Q10" = 1.45
    "R10" = 1.79
    "S10" = 2.16
    "T10" = 2.7
    "U10" = 1.52
    "V10" = 1.02
    "W10" = 8
    "Y10" = 110
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
